$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 150, pushing existing rows 150-152 down to 151-153
$ws.Rows.Item(150).Insert()

# Fill in the new row 150 with the new weekly record
$ws.Range("A150").Value = 1
$ws.Range("B150").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C150").Value = "Arica y Parinacota"
$ws.Range("D150").Value = 44448
$ws.Range("E150").Value = 15
$ws.Range("F150").Value = "Fruta"
$ws.Range("G150").Value = 100102
$ws.Range("H150").Value = "Cítricos"
$ws.Range("I150").Value = 100102003
$ws.Range("J150").Value = "Limón"
$ws.Range("K150").Value = "Sin especificar"
$ws.Range("L150").Value = "3a amarillo"
$ws.Range("M150").Value = 270
$ws.Range("N150").Value = 9500
$ws.Range("O150").Value = 10000
$ws.Range("P150").Value = 9750
$ws.Range("Q150").Value = "$/caja 20 kilos"
$ws.Range("R150").Value = "Región de Coquimbo"
$ws.Range("S150").Value = 488
$ws.Range("T150").Value = 20
